$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the week data between row 2 and row 3:
# Fecha (D), Volumen (J), Precio minimo (K), Precio maximo (L),
# Precio promedio ponderado (M), Precio $/Kg (P)

$row2Date = $ws.Range("D2").Value2
$row3Date = $ws.Range("D3").Value2
$ws.Range("D2").Value = $row3Date
$ws.Range("D3").Value = $row2Date

$row2J = $ws.Range("J2").Value2
$row3J = $ws.Range("J3").Value2
$ws.Range("J2").Value = $row3J
$ws.Range("J3").Value = $row2J

$row2K = $ws.Range("K2").Value2
$row3K = $ws.Range("K3").Value2
$ws.Range("K2").Value = $row3K
$ws.Range("K3").Value = $row2K

$row2L = $ws.Range("L2").Value2
$row3L = $ws.Range("L3").Value2
$ws.Range("L2").Value = $row3L
$ws.Range("L3").Value = $row2L

$row2M = $ws.Range("M2").Value2
$row3M = $ws.Range("M3").Value2
$ws.Range("M2").Value = $row3M
$ws.Range("M3").Value = $row2M

$row2P = $ws.Range("P2").Value2
$row3P = $ws.Range("P3").Value2
$ws.Range("P2").Value = $row3P
$ws.Range("P3").Value = $row2P
